# Add the new "Cotações atualizadas - 2025-10-11" row (row 37) to Sheet1,
# mirroring the existing rows' layout: column A is a serial date with the
# same number format/style as the row above it, columns B:E hold the
# quote values as text (decimal comma, Portuguese formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 37
$prevRow = $newRow - 1

# Date value for 2025-10-11 (Excel serial date number)
$ws.Cells.Item($newRow, 1).Value = 45941
# Match the date-formatted style used by the preceding date cell (A36)
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

# Quote values (kept as text, matching the existing inline-string cells)
$ws.Cells.Item($newRow, 2).Value = "21,6987"
$ws.Cells.Item($newRow, 3).Value = "15,6648"
$ws.Cells.Item($newRow, 4).Value = "15,4517"
$ws.Cells.Item($newRow, 5).Value = "15,4517"
